$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (kulbutoké)
$ws.Range("D3").Value = 33
$ws.Range("F3").Value = 15.93816666666667
$ws.Range("G3").Value = 734
$ws.Range("H3").Value = 362
$ws.Range("I3").Value = 67
$ws.Range("K3").Value = 4015
$ws.Range("L3").Value = 99
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 260
$ws.Range("P3").Value = 4.545454545454546
$ws.Range("Q3").Value = 7.878787878787879
$ws.Range("R3").Value = 22.24
$ws.Range("S3").Value = 28.98

# Row 4 (tomlora)
$ws.Range("D4").Value = 174
$ws.Range("E4").Value = 120
$ws.Range("F4").Value = 86.42216666666668
$ws.Range("G4").Value = 3615
$ws.Range("H4").Value = 1521
$ws.Range("I4").Value = 505
$ws.Range("J4").Value = 475
$ws.Range("K4").Value = 35171
$ws.Range("L4").Value = 1323
$ws.Range("M4").Value = 1004
$ws.Range("N4").Value = 1343
$ws.Range("O4").Value = 7.603448275862069
$ws.Range("P4").Value = 5.770114942528735
$ws.Range("Q4").Value = 7.718390804597701
$ws.Range("R4").Value = 20.78
$ws.Range("S4").Value = 29.8

# Row 5 (chatobogan)
$ws.Range("D5").Value = 54
$ws.Range("E5").Value = 45
$ws.Range("F5").Value = 28.504
$ws.Range("G5").Value = 3051
$ws.Range("H5").Value = 1468
$ws.Range("I5").Value = 422
$ws.Range("J5").Value = 390
$ws.Range("K5").Value = 5186
$ws.Range("L5").Value = 295
$ws.Range("M5").Value = 219
$ws.Range("N5").Value = 607
$ws.Range("O5").Value = 5.462962962962963
$ws.Range("P5").Value = 4.055555555555555
$ws.Range("Q5").Value = 11.24074074074074
$ws.Range("R5").Value = 56.5
$ws.Range("S5").Value = 31.67

# Row 6 (exorblue)
$ws.Range("D6").Value = 30
$ws.Range("F6").Value = 15.84316666666667
$ws.Range("G6").Value = 1490
$ws.Range("H6").Value = 561
$ws.Range("I6").Value = 199
$ws.Range("J6").Value = 215
$ws.Range("K6").Value = 3617
$ws.Range("L6").Value = 101
$ws.Range("M6").Value = 140
$ws.Range("N6").Value = 370
$ws.Range("O6").Value = 3.366666666666667
$ws.Range("P6").Value = 4.666666666666667
$ws.Range("Q6").Value = 12.33333333333333
$ws.Range("R6").Value = 49.67
$ws.Range("S6").Value = 31.69

# Row 9 (linò)
$ws.Range("D9").Value = 42
$ws.Range("E9").Value = 46
$ws.Range("F9").Value = 21.35016666666667
$ws.Range("G9").Value = 1422
$ws.Range("H9").Value = 483
$ws.Range("I9").Value = 201
$ws.Range("J9").Value = 122
$ws.Range("K9").Value = 9792
$ws.Range("L9").Value = 267
$ws.Range("M9").Value = 266
$ws.Range("N9").Value = 302
$ws.Range("O9").Value = 6.357142857142857
$ws.Range("P9").Value = 6.333333333333333
$ws.Range("Q9").Value = 7.190476190476191
$ws.Range("R9").Value = 33.86
$ws.Range("S9").Value = 30.5

# Row 10 (namiyeon)
$ws.Range("D10").Value = 78
$ws.Range("F10").Value = 38.55433333333333
$ws.Range("G10").Value = 6065
$ws.Range("H10").Value = 3494
$ws.Range("I10").Value = 773
$ws.Range("J10").Value = 823
$ws.Range("K10").Value = 1513
$ws.Range("L10").Value = 196
$ws.Range("M10").Value = 230
$ws.Range("N10").Value = 1242
$ws.Range("O10").Value = 2.512820512820513
$ws.Range("P10").Value = 2.948717948717949
$ws.Range("Q10").Value = 15.92307692307692
$ws.Range("R10").Value = 77.76000000000001
$ws.Range("S10").Value = 29.66

# Row 12 (chguizou)
$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 1.687666666666666
$ws.Range("G12").Value = 48
$ws.Range("H12").Value = 27
$ws.Range("I12").Value = 5
$ws.Range("K12").Value = 522
$ws.Range("L12").Value = 21
$ws.Range("M12").Value = 11
$ws.Range("N12").Value = 26
$ws.Range("O12").Value = 7
$ws.Range("P12").Value = 3.666666666666667
$ws.Range("Q12").Value = 8.666666666666666
$ws.Range("R12").Value = 16
$ws.Range("S12").Value = 33.75
